$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SearchData")
$ws.Range("A1").Value = "test"
